$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order: A=Expense Date, B=Expense Value, C=Description, D=Expense Type
$ws.Range("A1").Value = "Expense Date"
$ws.Range("B1").Value = "Expense Value"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Expense Type"

# New data row values
# A2 must remain a literal text string (not get auto-parsed into a date
# serial number) - enter it as a formula producing the literal text, then
# convert the cell to a plain value via copy / paste-special values so no
# extra number-format style gets attached to the cell.
$ws.Range("A2").Formula = "=""2024-03-14"""
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

$ws.Range("B2").Value = 123.0
$ws.Range("C2").Value = "q"
$ws.Range("D2").Value = "TRANSPORTATION"

# Column E is no longer used - clear it out entirely
$ws.Range("E1:E2").Clear()
